# Updates Tiamat_Profits market-price derived columns (H-N) for the rows
# that the scheduled price-refresh run touched, across the ALC/ARM/BSM/
# CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 5010.5
$ws.Range("I40").Value = 7928.7144
$ws.Range("J40").Value = 925
$ws.Range("K40").Value = 7928.7144
$ws.Range("L40").Value = 925
$ws.Range("M40").Value = -7753.7144
$ws.Range("N40").Value = -1275

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 26505.232
$ws.Range("I137").Value = 48478.383
$ws.Range("J137").Value = 5530.864
$ws.Range("K137").Value = 145435.149
$ws.Range("L137").Value = 16592.592
$ws.Range("M137").Value = -142885.149
$ws.Range("N137").Value = -21692.592

# Row 138: All-night Crafting
$ws.Range("H138").Value = 1611.2
$ws.Range("I138").Value = 843.3333
$ws.Range("J138").Value = 2167.2415
$ws.Range("K138").Value = 2529.9999
$ws.Range("L138").Value = 6501.7245
$ws.Range("M138").Value = 2610.0001
$ws.Range("N138").Value = -16781.7245

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 40953.19
$ws.Range("I74").Value = 47741.953
$ws.Range("J74").Value = 3615
$ws.Range("K74").Value = 47741.953
$ws.Range("L74").Value = 3615
$ws.Range("M74").Value = -46867.953
$ws.Range("N74").Value = -5363

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 40953.19
$ws.Range("I77").Value = 47741.953
$ws.Range("J77").Value = 3615
$ws.Range("K77").Value = 238709.765
$ws.Range("L77").Value = 18075
$ws.Range("M77").Value = -234341.765
$ws.Range("N77").Value = -26811

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1651036.5
$ws.Range("I132").Value = 1927494.9
$ws.Range("J132").Value = 674216.8
$ws.Range("K132").Value = 5782484.699999999
$ws.Range("L132").Value = 2022650.4
$ws.Range("M132").Value = -5779954.699999999
$ws.Range("N132").Value = -2027710.4

# Row 133: Shielding My Students
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 128: Mangalomania
$ws.Range("H128").Value = 1000
$ws.Range("I128").Value = 1000
$ws.Range("K128").Value = 3000
$ws.Range("M128").Value = -510

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 10995.214
$ws.Range("I31").Value = 11066.042
$ws.Range("J31").Value = 10900.777
$ws.Range("K31").Value = 11066.042
$ws.Range("L31").Value = 10900.777
$ws.Range("M31").Value = -10771.042
$ws.Range("N31").Value = -11490.777

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 10995.214
$ws.Range("I34").Value = 11066.042
$ws.Range("J34").Value = 10900.777
$ws.Range("K34").Value = 11066.042
$ws.Range("L34").Value = 10900.777
$ws.Range("M34").Value = -10864.042
$ws.Range("N34").Value = -11304.777

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 873.91113
$ws.Range("I134").Value = 732.1177
$ws.Range("J134").Value = 1312.1818
$ws.Range("K134").Value = 2196.3531
$ws.Range("L134").Value = 3936.5454
$ws.Range("M134").Value = 338.6468999999997
$ws.Range("N134").Value = -9006.545399999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 3229.1843
$ws.Range("I5").Value = 600.375
$ws.Range("J5").Value = 7735.7144
$ws.Range("K5").Value = 1801.125
$ws.Range("L5").Value = 23207.1432
$ws.Range("M5").Value = -1689.125
$ws.Range("N5").Value = -23431.1432

# Row 122: Salt of the North
$ws.Range("H122").Value = 313.88635
$ws.Range("I122").Value = 241.93103
$ws.Range("K122").Value = 2177.37927
$ws.Range("M122").Value = 272.6207300000001

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 89286230
$ws.Range("I131").Value = 401.125
$ws.Range("J131").Value = 208334000
$ws.Range("K131").Value = 1203.375
$ws.Range("L131").Value = 625002000
$ws.Range("M131").Value = 3836.625
$ws.Range("N131").Value = -625012080

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 3229.1843
$ws.Range("I135").Value = 600.375
$ws.Range("J135").Value = 7735.7144
$ws.Range("K135").Value = 5403.375
$ws.Range("L135").Value = 69621.4296
$ws.Range("M135").Value = -2868.375
$ws.Range("N135").Value = -74691.4296

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 1379.9592
$ws.Range("I126").Value = 1252.7222
$ws.Range("K126").Value = 3758.1666
$ws.Range("M126").Value = -1288.1666

$ws = $wb.Worksheets.Item("LTW")
# Row 9: From the Sands to the Stage
$ws.Range("H9").Value = 14758.5
$ws.Range("J9").Value = 35164.4
$ws.Range("L9").Value = 35164.4
$ws.Range("N9").Value = -35612.4

# Row 13: Throwing Down the Gauntlet
$ws.Range("H13").Value = 38108
$ws.Range("I13").Value = 5263
$ws.Range("J13").Value = 60004.668
$ws.Range("K13").Value = 5263
$ws.Range("L13").Value = 60004.668
$ws.Range("M13").Value = -5123
$ws.Range("N13").Value = -60284.668

# Row 20: Choke Hold
$ws.Range("H20").Value = 51486.332
$ws.Range("J20").Value = 51486.332
$ws.Range("L20").Value = 51486.332
$ws.Range("N20").Value = -51938.332

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1052.8667
$ws.Range("J22").Value = 1172.091
$ws.Range("L22").Value = 1172.091
$ws.Range("N22").Value = -1762.091

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1052.8667
$ws.Range("J27").Value = 1172.091
$ws.Range("L27").Value = 1172.091
$ws.Range("N27").Value = -1386.091

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 2399.3635
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 2732.5557
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 2732.5557
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -3108.5557

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 3311.111
$ws.Range("I61").Value = 5400
$ws.Range("J61").Value = 2714.2856
$ws.Range("K61").Value = 5400
$ws.Range("L61").Value = 2714.2856
$ws.Range("M61").Value = -5198
$ws.Range("N61").Value = -3118.2856

# Row 113: Peace in Rest
$ws.Range("H113").Value = 3311.111
$ws.Range("I113").Value = 5400
$ws.Range("J113").Value = 2714.2856
$ws.Range("K113").Value = 5400
$ws.Range("L113").Value = 2714.2856
$ws.Range("M113").Value = -3230
$ws.Range("N113").Value = -7054.2856

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 404960.66
$ws.Range("I132").Value = 137113.14
$ws.Range("J132").Value = 672808.2
$ws.Range("K132").Value = 411339.42
$ws.Range("L132").Value = 2018424.6
$ws.Range("M132").Value = -408809.42
$ws.Range("N132").Value = -2023484.6

$ws = $wb.Worksheets.Item("WVR")
# Row 30: The Telltale Tress
$ws.Range("H30").Value = 62006
$ws.Range("J30").Value = 62006
$ws.Range("L30").Value = 62006
$ws.Range("N30").Value = -62220

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 3393.2683
$ws.Range("I132").Value = 910.96
$ws.Range("J132").Value = 7271.875
$ws.Range("K132").Value = 2732.88
$ws.Range("L132").Value = 21815.625
$ws.Range("M132").Value = -202.8800000000001
$ws.Range("N132").Value = -26875.625

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 840206.25
$ws.Range("I136").Value = 1021103.2
$ws.Range("K136").Value = 3063309.6
$ws.Range("L136").Value = 1113642.66
$ws.Range("M136").Value = -3060759.6

Write-Output "Edits applied"
